# Natmi following Dr Hou advice
#
# Update the FAPs -> Rbp4/Stra6 -> FAPs / sCs edge rows with recomputed
# values, and add the corresponding new "M2" sending-cluster rows
# (M2 -> Rbp4/Stra6 -> FAPs / sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Rbp4/Stra6 -> FAPs (values recalculated)
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rbp4"
$ws.Range("C2").Value = "Stra6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.164095000000001
$ws.Range("H2").Value = 12.492285
$ws.Range("I2").Value = 0.9890011363911656
$ws.Range("J2").Value = 0.9890011363911656
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3437216666666667
$ws.Range("N2").Value = 1.031165
$ws.Range("O2").Value = 0.3023306815952884
$ws.Range("P2").Value = 0.3023306815952883
$ws.Range("Q2").Value = 1.431289673558334
$ws.Range("R2").Value = 12.881607062025
$ws.Range("S2").Value = 0.2990053876636559
$ws.Range("T2").Value = 0.2990053876636558

# Row 3: FAPs -> Rbp4/Stra6 -> sCs (values recalculated)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rbp4"
$ws.Range("C3").Value = "Stra6"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.164095000000001
$ws.Range("H3").Value = 12.492285
$ws.Range("I3").Value = 0.9890011363911656
$ws.Range("J3").Value = 0.9890011363911656
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7931846666666665
$ws.Range("N3").Value = 2.379554
$ws.Range("O3").Value = 0.6976693184047117
$ws.Range("P3").Value = 0.6976693184047117
$ws.Range("Q3").Value = 3.302896304543333
$ws.Range("R3").Value = 29.72606674089
$ws.Range("S3").Value = 0.6899957487275098
$ws.Range("T3").Value = 0.6899957487275098

# Row 4 (new): M2 -> Rbp4/Stra6 -> FAPs
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Rbp4"
$ws.Range("C4").Value = "Stra6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04630966666666667
$ws.Range("H4").Value = 0.138929
$ws.Range("I4").Value = 0.01099886360883443
$ws.Range("J4").Value = 0.01099886360883443
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3437216666666667
$ws.Range("N4").Value = 1.031165
$ws.Range("O4").Value = 0.3023306815952884
$ws.Range("P4").Value = 0.3023306815952883
$ws.Range("Q4").Value = 0.01591763580944444
$ws.Range("R4").Value = 0.143258722285
$ws.Range("S4").Value = 0.003325293931632526
$ws.Range("T4").Value = 0.003325293931632526

# Row 5 (new): M2 -> Rbp4/Stra6 -> sCs
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Rbp4"
$ws.Range("C5").Value = "Stra6"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04630966666666667
$ws.Range("H5").Value = 0.138929
$ws.Range("I5").Value = 0.01099886360883443
$ws.Range("J5").Value = 0.01099886360883443
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7931846666666665
$ws.Range("N5").Value = 2.379554
$ws.Range("O5").Value = 0.6976693184047117
$ws.Range("P5").Value = 0.6976693184047117
$ws.Range("Q5").Value = 0.03673211751844443
$ws.Range("R5").Value = 0.3305890576659999
$ws.Range("S5").Value = 0.007673569677201905
$ws.Range("T5").Value = 0.007673569677201906
